# Insert a new weekly price record as row 8 in the "Arveja Verde" price
# table, pushing the existing rows 8-34 down to rows 9-35 (dimension grows
# from A1:R34 to A1:R35).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 8:34 down one row to make room for the new record.
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the new weekly observation.
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C8").Value = "Ñuble"
$ws.Range("D8").Value2 = 44537
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 100112022
$ws.Range("G8").Value = "Arveja Verde"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 13000
$ws.Range("L8").Value = 14000
$ws.Range("M8").Value = 13500
$ws.Range("N8").Value = "$/saco 25 kilos"
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 540
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
